$wb = $excel.ActiveWorkbook

# --- JobType sheet: add 3 new "CF Job Type" rows (Strategy/Post Merger
# Integration/Valuation Advisory with their job codes) ---
$wsJobType = $wb.Worksheets.Item("JobType")

# Seed the JobCode column's format (B2 is stored as text via a custom
# number format) down onto the new rows first so "6291" etc. land as text,
# matching the existing "6421" entry, instead of turning into numbers.
$wsJobType.Range("B2").Copy($wsJobType.Range("B3"))
$wsJobType.Range("B2").Copy($wsJobType.Range("B4"))
$wsJobType.Range("B2").Copy($wsJobType.Range("B5"))

$wsJobType.Range("A3").Value = "Strategy"
$wsJobType.Range("C3").Value = "CF Job Type"
$wsJobType.Range("B3").Value = "6291"

$wsJobType.Range("A4").Value = "Post Merger Integration"
$wsJobType.Range("B4").Value = "6301"

$wsJobType.Range("A5").Value = "Valuation Advisory"
$wsJobType.Range("B5").Value = "6311"

$wsJobType.Range("C4").Value = "CF Job Type"
$wsJobType.Range("C5").Value = "CF Job Type"

$wsJobType.Range("A3:C3").RowHeight = 15.75
$wsJobType.Range("A4:C4").RowHeight = 15.75
$wsJobType.Range("A5:C5").RowHeight = 15.75

# --- Users sheet: add 3 new rows (Matthew Gemmell / Meissa Lee) ---
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("B3").Value = "Meissa Lee"
$wsUsers.Range("A3").Value = "Matthew Gemmell"
$wsUsers.Range("A4").Value = "Matthew Gemmell"
$wsUsers.Range("B4").Value = "Meissa Lee"
$wsUsers.Range("A5").Value = "Matthew Gemmell"
$wsUsers.Range("B5").Value = "Meissa Lee"

$wsUsers.Columns.Item(1).AutoFit()

# --- Selection / active sheet changes ---
$wsJobType.Range("A2:A5").Select()
$wsUsers.Range("A2:A5").Select()
$wsUsers.Activate()
$wsUsers.Select()
